$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, alphabetically-resorted list of cluster names (row 2 downwards) together
# with their updated "Active cases" counts. Three old clusters were dropped
# ("126 Racecourse Road Public Housing Tower Flemington", "FedEx Station
# Melbourne Airport", "Nido Early School Glenroy") and two new ones were
# added ("Hello Fresh Warehouse Ravenhall", "Launch Housing City Edge Crisis
# Accommodation South Melbourne"), and every remaining count was refreshed.
$names = @(
    "139 Highett St Apartment Complex Richmond",
    "3175 The Bays Aged Care Facility Hastings",
    "3600 Belvedere Age Care Noble Park Outbreak",
    "3612 BlueCross Glengowrie Outbreak",
    "3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak",
    "4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak",
    "Australian Lamb Colac East",
    "Bread Solutions Braeside Outbreak",
    "CS Square Caroline Springs Outbreak",
    "Cedar Meats Australia Brooklyn Outbreak",
    "Child's Play Early Learning Centre Tarneit",
    "Embracia Aged Care Reservoir Outbreak",
    "Guardian Childcare Caulfield Outbreak",
    "Hello Fresh Warehouse Ravenhall",
    "Inghams Enterprise Somerville Outbreak",
    "Kool Kidz Childcare Narre Warren",
    "Lantmannen Unibake Australia Mordialloc",
    "Launch Housing City Edge Crisis Accommodation South Melbourne",
    "Nido Early School Ascot Vale",
    "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
    "Northern Health The Northern Hospital Epping",
    "Oceania Meat Processors Laverton North Outbreak",
    "Pick It Up Fitness Mulgrave Outbreak",
    "Robin Hood Inn Drouin West Outbreak",
    "Social Gathering Warrnambool 28 Sep Outbreak",
    "St Vincents Hospital Emergency Department Melbourne",
    "Target Distribution Centre Truganina Outbreak",
    "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B",
    "The Toolshed Bar Private Event Noojee",
    "Turosi Breakwater",
    "Visy Recycling Springvale",
    "Werribee Mercy Hospital Emergency Department",
    "Western Health Sunshine Hospital Emergency Department"
)

$values = @(11,16,21,30,20,16,13,19,17,11,11,22,17,13,15,10,25,12,11,42,15,16,11,41,13,42,20,12,17,10,29,24,21)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# The sheet used to have 34 data rows (row 35 was the last one); now it only
# has 33, so drop the now-unused trailing row entirely.
$ws.Range("A35:B35").ClearContents()
